# Applies the cryptos.xlsx price/volume/ranking refresh described by the
# commit "Updated cryptos list on Mon Jan  1 11:12:08 UTC 2024 with GitHub Actions".
#
# Every data row (spreadsheet rows 2-51) gets its Price (column D) and
# Volume(1h) (column E) refreshed with the latest scrape. Two coin pairs also
# swapped ranking order in this refresh (InternetComputer(DFINITY)/ShibaInu at
# rows 20-21, and FraxShare/THORChain at rows 50-51), so their Coin name
# (column B) and Link (column C) cells are rewritten too, in addition to D/E.
#
# Price strings that look like a plain decimal number (a single '.') are
# written with a leading apostrophe so Excel stores them as literal text --
# matching the source workbook, where every Price cell is text (e.g.
# "310.40", "0.0000105") -- instead of silently coercing the cell to a
# Number. Prices already punctuated with thousands-separator dots (e.g.
# "42.762.71") can never parse as a number, so they don't need the guard.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.762.71"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.304.31"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  -0.40%  "
$ws.Range("D5").Value = "'310.40"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").Value = "'104.88"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  -0.83%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "'39.74"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "'8.28"
$ws.Range("E12").Value = "  -3.67%  "
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").Value = "'0.994"
$ws.Range("E14").Value = "  +1.76%  "
$ws.Range("D15").Value = "2.782.60"
$ws.Range("E15").Value = "  +4.48%  "
$ws.Range("D16").Value = "'15.38"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "2.303.07"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "42.702.14"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "'7.34"
$ws.Range("E19").Value = "  -2.68%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'13.66"
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0000105"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").Value = "'73.44"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").Value = "'3.44"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("D24").Value = "'268.11"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "'7.24"
$ws.Range("E28").Value = "  +15.37%  "
$ws.Range("D29").Value = "'2.25"
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("D31").Value = "'36.21"
$ws.Range("E31").Value = "  -5.16%  "
$ws.Range("D32").Value = "'164.81"
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").Value = "'0.0859"
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("E35").Value = "  +0.87%  "
$ws.Range("E36").Value = "  -3.91%  "
$ws.Range("D37").Value = "'4.54"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("D38").Value = "'0.0349"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").Value = "'2.83"
$ws.Range("E39").Value = "  +3.80%  "
$ws.Range("D40").Value = "'3.62"
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("D41").Value = "'110.44"
$ws.Range("E41").Value = "  +12.14%  "
$ws.Range("D42").Value = "'1.59"
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("D43").Value = "'70.93"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'12.26"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").Value = "1.723.96"
$ws.Range("E47").Value = "  +6.46%  "
$ws.Range("D48").Value = "'110.89"
$ws.Range("E48").Value = "  -3.37%  "
$ws.Range("D49").Value = "'77.86"
$ws.Range("E49").Value = "  -5.43%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'8.66"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'5.14"
$ws.Range("E51").Value = "  -3.33%  "
